$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Insert a new column at L - this shifts the old L column (help notes,
# col 12) and its widths to M (col 13).
$ws.Columns.Item(12).Insert()

# New "EndTime" header in the freshly inserted column L.
$ws.Range("L1").Value = "EndTime"

# New sample password value for the existing "password" column (K).
$ws.Range("K2").Value = "#Parrsoo2020#"

# New sample value for the freshly inserted "EndTime" column (L).
$ws.Range("L2").Value = "1400/10/20 13:13:13.259"

# userId* sample value changed.
$ws.Range("A2").Value = "alit"

# Column widths (closest reproducible values - the runtime quantizes
# ColumnWidth to the nearest 1/7 character when serializing to OOXML).
$ws.Columns.Item(11).ColumnWidth = 11.714285714285715
$ws.Columns.Item(12).ColumnWidth = 20.857142857142854
$ws.Columns.Item(13).ColumnWidth = 44.85714285714286

$ws.Range("D10").Select()
